# Append two new "ListParagraph" bullet items after the existing
# "Переделать USB: ..." item, and carry the _GoBack bookmark along so it
# ends up collapsed at the very end of the new last paragraph's text
# (exactly as it originally sat at the end of the old last paragraph).

$d = $word.ActiveDocument

$text1 = "Пересмотреть питание и зарядку. А то зарядка может не идти, если прочее будет потреблять много."
$text2 = "Переделать инфракрасный светодиод."
$cr = [char]13

# The _GoBack bookmark is currently a collapsed range sitting right after
# the text of the last paragraph (before its paragraph mark). Grab a
# Range at that exact spot so we can grow new content "in front of" it;
# Word keeps a collapsed bookmark glued to the text that follows it, so
# text inserted via InsertBefore() at this point ends up before the
# bookmark and the bookmark rides along to stay at the end.
$bm = $d.Bookmarks.Item("_GoBack")
$anchor = $d.Range($bm.Start, $bm.Start)

# First chunk: a paragraph break (cloning the current ListParagraph /
# numbering / rPr), the first new bullet's text, and another paragraph
# break to start the second bullet. Keeping this to a single
# InsertBefore call correctly reproduces the paragraph formatting
# (style + numPr + rPr) on both new paragraphs.
$anchor.InsertBefore($cr + $text1 + $cr)

# Re-collapse to the (now shifted) original anchor point - this is where
# the second bullet's text must land, immediately in front of the
# bookmark.
$anchor.Collapse(0)
$anchor.InsertBefore($text2)

# The run that now holds $text2 needs the same Russian-language direct
# formatting as its sibling runs; set it explicitly (InsertBefore after
# a fresh Collapse() does not always inherit the preceding run's rPr).
$anchor.LanguageID = "ru-RU"
